# tests/reports/df_filter_frame.xlsx
#
# The sheet used to carry a standalone "<frame>" marker cell in A1 (row 1)
# of Sheet1. This edit turns that marker into a cell Comment/Note attached
# to what becomes the new A1 (the old "{{ mystring }}" cell), and removes
# the now-empty marker row, shifting every following row up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Drop the old row 2 (the "<frame>" marker lived alone in row 1; row 2 held
# "{{ mystring }}"). Deleting row 2 shifts rows 4,5,7,9 up to 3,4,6,8 and
# keeps A1's original style/text untouched for now.
$ws.Rows.Item(2).Delete() | Out-Null

# Re-purpose A1: attach the frame marker as a cell comment/note, then
# overwrite the cell value with what used to live in row 2.
$ws.Range("A1").AddComment("<frame>") | Out-Null
$ws.Range("A1").Value = "{{ mystring }}"

# Restore the (otherwise lost) active-cell selection on Sheet1.
$ws.Range("B13").Select() | Out-Null
